$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 187
$ws.Range("E2").Value = 94
$ws.Range("I2").Value = 9051688
$ws.Range("J2").Value = 3909033
$ws.Range("N2").Value = 69.54000000000001
$ws.Range("O2").Value = 30.03
$ws.Range("S2").Value = 30.03
$ws.Range("U2").Value = 1794925

# Row 3 updates
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 2
$ws.Range("L3").Value = 353390
$ws.Range("M3").Value = 119941
$ws.Range("Q3").Value = 7.65
$ws.Range("R3").Value = 2.6
$ws.Range("S3").Value = 7.65
$ws.Range("T3").Value = 7.65
$ws.Range("U3").Value = 283700
